$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "23 x 21" + [char]11 + "  2    1" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "3|    |"
$t.Cell(1,2).Range.Text = "46 x 56" + [char]11 + "  5    6" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "6|    |"
$t.Cell(1,3).Range.Text = "77 x 75" + [char]11 + "  7    5" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "7|    |"
$t.Cell(2,1).Range.Text = "59 x 68" + [char]11 + "  6    8" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "9|    |"
$t.Cell(2,2).Range.Text = "59 x 37" + [char]11 + "  3    7" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "9|    |"
$t.Cell(2,3).Range.Text = "62 x 11" + [char]11 + "  1    1" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "2|    |"
$t.Cell(3,1).Range.Text = "44 x 94" + [char]11 + "  9    4" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "4|    |"
$t.Cell(3,2).Range.Text = "77 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "7|    |"
$t.Cell(3,3).Range.Text = "39 x 55" + [char]11 + "  5    5" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "9|    |"
$t.Cell(4,1).Range.Text = "97 x 10" + [char]11 + "  1    0" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(4,2).Range.Text = "79 x 26" + [char]11 + "  2    6" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "9|    |"
$t.Cell(4,3).Range.Text = "35 x 59" + [char]11 + "  5    9" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"
$t.Cell(5,1).Range.Text = "88 x 56" + [char]11 + "  5    6" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "8|    |"
$t.Cell(5,2).Range.Text = "53 x 74" + [char]11 + "  7    4" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "3|    |"
$t.Cell(5,3).Range.Text = "15 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "5|    |"
